# Actualización automática 2025-09-30 17:30:09
#
# Updates the sales figures for ALMEIDA CUATIN JHONATHANN CARLOS /
# MANCHENO PINO HERVIN SANTIAGO on the "VENTAS POR GRUPO" sheet, and the
# derived monthly / compliance totals on the other two sheets (these are
# plain cached values in the workbook, not live formulas, so each
# dependent cell is updated explicitly to keep the workbook internally
# consistent).

$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO": raw sales figures for row 20 -----------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("E20").Value = 220.24
$wsGrupo.Range("H20").Value = 205.2

# --- Sheet "VENTA MENSUAL": septiembre column for the same client -----
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F20").Value = 4354.56
$wsMensual.Range("F35").Value = 29882.18

# --- Sheet "CUMPLIMIENTO MENSUAL": recomputed compliance figures ------
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Row 4 - FREGADEROS DE COCINA
$wsCumplimiento.Range("D4").Value = 363.66
$wsCumplimiento.Range("E4").Value = 157.95144263264
$wsCumplimiento.Range("F4").Value = 0.6971856256921076

# Row 6 - INODOROS
$wsCumplimiento.Range("D6").Value = 1010.7
$wsCumplimiento.Range("E6").Value = -196.5765691911271
$wsCumplimiento.Range("F6").Value = 1.24145794329469

# Row 15 - TOTAL
$wsCumplimiento.Range("D15").Value = 30139.41
$wsCumplimiento.Range("E15").Value = 8603.608813395924
$wsCumplimiento.Range("F15").Value = 0.7779313776545178
